$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Infused Beverages")

# --- Insert the "sawblade_holly_berry_cider" row right after the cider rows (new row 53) ---
$ws.Rows("53:53").Insert()

$ws.Range("A53").Value2 = "sawblade_holly_berry_cider"
$ws.Range("B53").Value2 = "toughness"
$ws.Range("C53").Value2 = "strength"
$ws.Range("F53").Value2 = "weakness"
$ws.Range("G53").Value2 = "nausea"
$ws.Range("I53").Value2 = "sawblade_holly_berry"

# rename glass_peach_cider's second effect from toughness -> projectile_rebound
$ws.Range("B52").Value2 = "projectile_rebound"

# --- Insert the "sawblade_holly_berry_liquor" row right after the liquor rows (new row 61) ---
$ws.Rows("61:61").Insert()

$ws.Range("A61").Value2 = "sawblade_holly_berry_liquor"
$ws.Range("B61").Value2 = "toughness"
$ws.Range("C61").Value2 = "speed"
$ws.Range("F61").Value2 = "slowness"
$ws.Range("G61").Value2 = "nausea"
$ws.Range("I61").Value2 = "sawblade_holly_berry, sugar"

# rename glass_peach_liquor's second effect from toughness -> projectile_rebound
$ws.Range("B60").Value2 = "projectile_rebound"

$ws.Range("I61").Select()
